$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("G4").Value = 2.9
$ws.Range("H4").Value = 2.7
$ws.Range("I4").Value = 2.72
$ws.Range("J4").Value = 3.55
$ws.Range("L4").Value = 3.3
$ws.Range("M4").Value = 1.12
$ws.Range("N4").Value = 5.4
$ws.Range("O4").Value = 1.47
$ws.Range("P4").Value = 2.5
$ws.Range("Q4").Value = 2.4
$ws.Range("R4").Value = 1.5
$ws.Range("U4").Value = 4.2
$ws.Range("V4").Value = 1.19
$ws.Range("Y4").Value = 1.93
$ws.Range("Z4").Value = 1.78
$ws.Range("AA4").Value = 7
$ws.Range("AB4").Value = 13.5
$ws.Range("AC4").Value = 10.75
$ws.Range("AD4").Value = 37
$ws.Range("AE4").Value = 30
$ws.Range("AF4").Value = 45
$ws.Range("AG4").Value = 5.4
$ws.Range("AH4").Value = 5.3
$ws.Range("AI4").Value = 15
$ws.Range("AJ4").Value = 90
$ws.Range("AK4").Value = 800
$ws.Range("AL4").Value = 7
$ws.Range("AM4").Value = 13
$ws.Range("AN4").Value = 10
$ws.Range("AO4").Value = 35
$ws.Range("AP4").Value = 26
$ws.Range("AQ4").Value = 40

# Row 5
$ws.Range("G5").Value = 3.45
$ws.Range("H5").Value = 2.7
$ws.Range("I5").Value = 2.37
$ws.Range("J5").Value = 4.15
$ws.Range("K5").Value = 1.82
$ws.Range("L5").Value = 3.1
$ws.Range("M5").Value = 1.14
$ws.Range("N5").Value = 4.9
$ws.Range("O5").Value = 1.57
$ws.Range("P5").Value = 2.25
$ws.Range("Q5").Value = 2.67
$ws.Range("R5").Value = 1.42
$ws.Range("U5").Value = 4.85
$ws.Range("V5").Value = 1.14
$ws.Range("W5").Value = 1.6
$ws.Range("X5").Value = 2.18
$ws.Range("Y5").Value = 2.15
$ws.Range("Z5").Value = 1.62
$ws.Range("AB5").Value = 16.5
$ws.Range("AC5").Value = 12.5
$ws.Range("AD5").Value = 50
$ws.Range("AE5").Value = 40
$ws.Range("AF5").Value = 60
$ws.Range("AG5").Value = 4.9
$ws.Range("AH5").Value = 5.5
$ws.Range("AI5").Value = 18
$ws.Range("AJ5").Value = 120
$ws.Range("AL5").Value = 5.7
$ws.Range("AM5").Value = 10
$ws.Range("AN5").Value = 9.75
$ws.Range("AO5").Value = 26
$ws.Range("AP5").Value = 25

# Row 6
$ws.Range("G6").Value = 1.82
$ws.Range("H6").Value = 3.15
$ws.Range("I6").Value = 4.7
$ws.Range("J6").Value = 2.42
$ws.Range("L6").Value = 4.9
$ws.Range("N6").Value = 6.2
$ws.Range("O6").Value = 1.39
$ws.Range("P6").Value = 2.77
$ws.Range("U6").Value = 3.65
$ws.Range("V6").Value = 1.24
$ws.Range("W6").Value = 1.44
$ws.Range("X6").Value = 2.57
$ws.Range("Y6").Value = 1.93
$ws.Range("Z6").Value = 1.78
$ws.Range("AA6").Value = 5.8
$ws.Range("AB6").Value = 7.7
$ws.Range("AD6").Value = 15
$ws.Range("AE6").Value = 16
$ws.Range("AG6").Value = 6.2
$ws.Range("AH6").Value = 6.1
$ws.Range("AI6").Value = 16
$ws.Range("AJ6").Value = 80
$ws.Range("AL6").Value = 11.5
$ws.Range("AM6").Value = 27
$ws.Range("AN6").Value = 15
$ws.Range("AO6").Value = 90
$ws.Range("AP6").Value = 50

# Row 7
$ws.Range("G7").Value = 1.65
$ws.Range("H7").Value = 3.7
$ws.Range("I7").Value = 5.25
$ws.Range("J7").Value = 2.3
$ws.Range("K7").Value = 2.1
$ws.Range("L7").Value = 6
$ws.Range("M7").Value = 1.07
$ws.Range("N7").Value = 9
$ws.Range("O7").Value = 1.36
$ws.Range("P7").Value = 3.2
$ws.Range("Q7").Value = 2.2
$ws.Range("R7").Value = 1.67
$ws.Range("U7").Value = 4
$ws.Range("V7").Value = 1.25
$ws.Range("W7").Value = 1.44
$ws.Range("X7").Value = 2.63
$ws.Range("AA7").Value = 6
$ws.Range("AB7").Value = 7
$ws.Range("AC7").Value = 8.5
$ws.Range("AD7").Value = 12
$ws.Range("AE7").Value = 15
$ws.Range("AG7").Value = 8.5
$ws.Range("AH7").Value = 7
$ws.Range("AL7").Value = 12
$ws.Range("AM7").Value = 26

# Row 8
$ws.Range("O8").Value = 1.29
$ws.Range("P8").Value = 3.75
$ws.Range("Q8").Value = 1.9
$ws.Range("R8").Value = 1.95
$ws.Range("W8").Value = 1.4
$ws.Range("Y8").Value = 3.25
$ws.Range("Z8").Value = 1.33
$ws.Range("AA8").Value = 5
$ws.Range("AC8").Value = 11
$ws.Range("AE8").Value = 15
$ws.Range("AG8").Value = 9.5
$ws.Range("AH8").Value = 13
$ws.Range("AJ8").Value = 201

# Row 9
$ws.Range("G9").Value = 1.57
$ws.Range("H9").Value = 4
$ws.Range("K9").Value = 2.2
$ws.Range("M9").Value = 1.06
$ws.Range("N9").Value = 9.5
$ws.Range("O9").Value = 1.33
$ws.Range("P9").Value = 3.4
$ws.Range("Q9").Value = 2
$ws.Range("R9").Value = 1.8
$ws.Range("U9").Value = 3.75
$ws.Range("V9").Value = 1.29
$ws.Range("W9").Value = 1.44
$ws.Range("X9").Value = 2.63
$ws.Range("AB9").Value = 7
$ws.Range("AE9").Value = 13
$ws.Range("AF9").Value = 29
$ws.Range("AG9").Value = 9.5
$ws.Range("AK9").Value = 451
$ws.Range("AL9").Value = 13
$ws.Range("AP9").Value = 41

# Row 10
$ws.Range("G10").Value = 2.35
$ws.Range("H10").Value = 3.2
$ws.Range("J10").Value = 3.1
$ws.Range("K10").Value = 2
$ws.Range("M10").Value = 1.08
$ws.Range("N10").Value = 8
$ws.Range("O10").Value = 1.4
$ws.Range("P10").Value = 3
$ws.Range("Q10").Value = 2.2
$ws.Range("R10").Value = 1.67
$ws.Range("U10").Value = 4.33
$ws.Range("V10").Value = 1.22
$ws.Range("W10").Value = 1.5
$ws.Range("X10").Value = 2.5
$ws.Range("Y10").Value = 1.83
$ws.Range("Z10").Value = 1.83
$ws.Range("AA10").Value = 7
$ws.Range("AD10").Value = 23
$ws.Range("AE10").Value = 21
$ws.Range("AF10").Value = 34
$ws.Range("AG10").Value = 8
$ws.Range("AH10").Value = 6
$ws.Range("AK10").Value = 351
$ws.Range("AL10").Value = 8.5
$ws.Range("AQ10").Value = 41

# Row 11
$ws.Range("Q11").Value = 2.3
$ws.Range("R11").Value = 1.62
$ws.Range("U11").Value = 4.33
$ws.Range("V11").Value = 1.22

# Row 12
$ws.Range("G12").Value = 2
$ws.Range("I12").Value = 3.8
$ws.Range("Q12").Value = 1.9
$ws.Range("R12").Value = 1.95
$ws.Range("AG12").Value = 11
$ws.Range("AI12").Value = 12
$ws.Range("AL12").Value = 12

# Row 13
$ws.Range("I13").Value = 3.3
$ws.Range("Q13").Value = 1.75
$ws.Range("R13").Value = 2.05
$ws.Range("W13").Value = 1.36
$ws.Range("X13").Value = 3
$ws.Range("Y13").Value = 1.67
$ws.Range("Z13").Value = 2.1
$ws.Range("AA13").Value = 8.5

# Row 14
$ws.Range("G14").Value = 1.44
$ws.Range("H14").Value = 4.75
$ws.Range("I14").Value = 5.5
$ws.Range("J14").Value = 1.95
$ws.Range("K14").Value = 2.6
$ws.Range("L14").Value = 5.5
$ws.Range("S14").Value = 1.8
$ws.Range("T14").Value = 2.05
$ws.Range("U14").Value = 2.2
$ws.Range("V14").Value = 1.62
$ws.Range("Y14").Value = 1.67
$ws.Range("Z14").Value = 2.1
$ws.Range("AA14").Value = 9.5
$ws.Range("AB14").Value = 8.5
$ws.Range("AD14").Value = 11
$ws.Range("AI14").Value = 17
$ws.Range("AL14").Value = 21
$ws.Range("AM14").Value = 34
$ws.Range("AO14").Value = 67
$ws.Range("AP14").Value = 41
$ws.Range("AQ14").Value = 41

# Row 15
$ws.Range("G15").Value = 2.7
$ws.Range("H15").Value = 3.1
$ws.Range("I15").Value = 2.6
$ws.Range("J15").Value = 3.25
$ws.Range("K15").Value = 2.2
$ws.Range("L15").Value = 3.1
$ws.Range("M15").Value = 1.05
$ws.Range("N15").Value = 11
$ws.Range("O15").Value = 1.25
$ws.Range("P15").Value = 3.75
$ws.Range("Q15").Value = 1.85
$ws.Range("R15").Value = 2
$ws.Range("U15").Value = 3
$ws.Range("V15").Value = 1.36
$ws.Range("W15").Value = 1.36
$ws.Range("X15").Value = 3
$ws.Range("Y15").Value = 1.62
$ws.Range("Z15").Value = 2.2
$ws.Range("AC15").Value = 11
$ws.Range("AE15").Value = 21
$ws.Range("AG15").Value = 11
$ws.Range("AH15").Value = 6
$ws.Range("AI15").Value = 12
$ws.Range("AK15").Value = 151
$ws.Range("AL15").Value = 10
$ws.Range("AM15").Value = 13

# Row 16
$ws.Range("G16").Value = 2.25
$ws.Range("H16").Value = 3.3
$ws.Range("I16").Value = 2.9
$ws.Range("J16").Value = 3
$ws.Range("L16").Value = 3.6
$ws.Range("U16").Value = 3.75
$ws.Range("V16").Value = 1.25
$ws.Range("Y16").Value = 1.8
$ws.Range("Z16").Value = 1.91
$ws.Range("AB16").Value = 11
$ws.Range("AO16").Value = 29
$ws.Range("AP16").Value = 23

